# ---------------------------------------------------------------------------
# Updated Columnar output from SAS with better formats
#
# This script reproduces, via the Excel object model, the changes made to
# test-columnar.xlsx:
#   1. The "testing" defined name now quotes the sheet name
#      (testing!$A$1:$H$12 -> 'testing'!$A$1:$H$12).
#   2. The workbook is flagged to fully recalculate the next time it is
#      opened (mirrors <calcPr fullCalcOnLoad="true"/> replacing the old
#      cached calcId).
#   3. The Percent Done / Percent Increment / datetime / time columns
#      (D, E, G, H) are rewritten with their full double-precision values
#      instead of the previous ~7 significant digit rounded values, as
#      produced by the newer/better SAS export formats.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Quote the sheet name inside the "testing" defined name's reference.
$nm = $wb.Names.Item("testing")
$nm.RefersTo = "='testing'!" + '$A$1:$H$12'

# 2) Request a full recalculation the next time the workbook is opened.
$wb.ForceFullCalculation = $true

# 3) Rewrite the Percent Done (D), Percent Increment (E), datetime (G) and
#    time (H) columns with full double precision values.
    # Row 2: G2=21916.000011574073, H2=0.000011574074074074073
    $ws.Range("G2").Value2 = 21916.000011574073
    $ws.Range("H2").Value2 = 0.000011574074074074073
    # Row 3: D3=0.1, G3=21916.000115740742, H3=0.00003472222222222222
    $ws.Range("D3").Value2 = 0.1
    $ws.Range("G3").Value2 = 21916.000115740742
    $ws.Range("H3").Value2 = 0.00003472222222222222
    # Row 4: D4=0.2, E4=0.5, G4=21916.00115740741, H4=0.00010416666666666666
    $ws.Range("D4").Value2 = 0.2
    $ws.Range("E4").Value2 = 0.5
    $ws.Range("G4").Value2 = 21916.00115740741
    $ws.Range("H4").Value2 = 0.00010416666666666666
    # Row 5: D5=0.3, E5=0.6666666666666666, G5=21916.011574074077, H5=0.0003125
    $ws.Range("D5").Value2 = 0.3
    $ws.Range("E5").Value2 = 0.6666666666666666
    $ws.Range("G5").Value2 = 21916.011574074077
    $ws.Range("H5").Value2 = 0.0003125
    # Row 6: D6=0.4, E6=0.75, G6=21916.11574074074, H6=0.0009375000000000001
    $ws.Range("D6").Value2 = 0.4
    $ws.Range("E6").Value2 = 0.75
    $ws.Range("G6").Value2 = 21916.11574074074
    $ws.Range("H6").Value2 = 0.0009375000000000001
    # Row 7: D7=0.5, E7=0.8, G7=21917.15740740741, H7=0.0028125
    $ws.Range("D7").Value2 = 0.5
    $ws.Range("E7").Value2 = 0.8
    $ws.Range("G7").Value2 = 21917.15740740741
    $ws.Range("H7").Value2 = 0.0028125
    # Row 8: D8=0.6, E8=0.8333333333333334, G8=21927.574074074077, H8=0.0084375
    $ws.Range("D8").Value2 = 0.6
    $ws.Range("E8").Value2 = 0.8333333333333334
    $ws.Range("G8").Value2 = 21927.574074074077
    $ws.Range("H8").Value2 = 0.0084375
    # Row 9: D9=0.7, E9=0.8571428571428571, G9=22031.74074074074, H9=0.0253125
    $ws.Range("D9").Value2 = 0.7
    $ws.Range("E9").Value2 = 0.8571428571428571
    $ws.Range("G9").Value2 = 22031.74074074074
    $ws.Range("H9").Value2 = 0.0253125
    # Row 10: D10=0.8, E10=0.875, G10=23073.40740740741, H10=0.07593749999999999
    $ws.Range("D10").Value2 = 0.8
    $ws.Range("E10").Value2 = 0.875
    $ws.Range("G10").Value2 = 23073.40740740741
    $ws.Range("H10").Value2 = 0.07593749999999999
    # Row 11: D11=0.9, E11=0.8888888888888888, G11=33490.07407407407, H11=0.2278125
    $ws.Range("D11").Value2 = 0.9
    $ws.Range("E11").Value2 = 0.8888888888888888
    $ws.Range("G11").Value2 = 33490.07407407407
    $ws.Range("H11").Value2 = 0.2278125
    # Row 12: E12=0.9, G12=137656.74074074073, H12=0.6834375
    $ws.Range("E12").Value2 = 0.9
    $ws.Range("G12").Value2 = 137656.74074074073
    $ws.Range("H12").Value2 = 0.6834375


$wb.Saved = $false
